$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo "testLadtNameN" -> "testLastNameN" in column L (rows 1-9)
for ($i = 1; $i -le 9; $i++) {
    $ws.Range("L$i").Value = "testLastName$i"
}

# Update the selected range to reflect column L instead of column M
$ws.Range("L1:L9").Select()
